$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 4 (shifts old rows 4..29 down to 6..31),
# pushing the existing "HKL" series down by two positions.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# The Insert() above mints a border-less style for the new rows' column A
# cells; copy the (bold/border/centered) formatting used by the rest of
# column A so the new index cells look like all the others.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 4: index 2 -> "Holden"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.04469962601295
$ws.Range("D4").Value = 0.9155068659708696
$ws.Range("E4").Value = 0.8511668408555053
$ws.Range("F4").Value = 0.8511668408555053
$ws.Range("G4").Value = 0.9233904661497294
$ws.Range("H4").Value = 0.9136629262329522
$ws.Range("I4").Value = 1.200564689223873
$ws.Range("J4").Value = 1.063831493884743
$ws.Range("K4").Value = 0.9808824174017727
$ws.Range("L4").Value = 1.044699626013105
$ws.Range("M4").Value = 1.044699626013105
$ws.Range("N4").Value = 1.044699626013105
$ws.Range("O4").Value = 0.9155068659708696
$ws.Range("P4").Value = 0.8833368534131875
$ws.Range("Q4").Value = 0.9896691799278061
$ws.Range("R4").Value = 0.9371244442798264
$ws.Range("S4").Value = 0.9435017335703725
$ws.Range("T4").Value = 0.9371244442798264
$ws.Range("U4").Value = 0.9688012066810554
$ws.Range("V4").Value = 0.9839808905474652
$ws.Range("W4").Value = 0.9867131657165686

# New row 5: index 3 -> "Rizzie Spiral"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 2.830978484609258
$ws.Range("D5").Value = 0.9549311586467293
$ws.Range("E5").Value = 1.893381439646856
$ws.Range("F5").Value = 1.893381439646856
$ws.Range("G5").Value = 0.2631854167387073
$ws.Range("H5").Value = 1.096876578893782
$ws.Range("I5").Value = 1.043346781155085
$ws.Range("J5").Value = 0.8810825262641679
$ws.Range("K5").Value = 1.08844351561893
$ws.Range("L5").Value = 2.830978484609258
$ws.Range("M5").Value = 2.830978484609258
$ws.Range("N5").Value = 2.830978484609258
$ws.Range("O5").Value = 0.9549311586467293
$ws.Range("P5").Value = 1.424156299146793
$ws.Range("Q5").Value = 0.9180068424554486
$ws.Range("R5").Value = 1.893097027634281
$ws.Range("S5").Value = 1.243131708185918
$ws.Range("T5").Value = 1.893097027634282
$ws.Range("U5").Value = 1.640093402291753
$ws.Range("V5").Value = 1.878270418755254
$ws.Range("W5").Value = 1.25652823769669

# Rename the "Thomas Hex" entry (now on row 11, after the insert shift) to
# "Matthies Hex".
$ws.Range("B11").Value = "Matthies Hex"
